$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-point the three tables (slides 14, 15, 16) at the built-in
#    "No Style, Table Grid" table style instead of the custom Google-Slides
#    generated style.
# ---------------------------------------------------------------------------
$newStyleId = "{E4DE9A76-1337-4772-96BF-B7339E8F785B}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's colour scheme back from the applied "Integral" design
#    (Red Violet colours) to the default "Office Theme" colours. The colour
#    scheme is shared by every slide/layout off the single slide master, so
#    editing it once (through any slide) updates the whole presentation.
# ---------------------------------------------------------------------------
# COM order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
